$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 (I0) and J1 (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style/formatting (bold font, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the I0 / IF data values for rows 2-21
$data = @{
    2  = @(7, 9)
    3  = @(7, 8)
    4  = @(7, 8)
    5  = @(5, 8)
    6  = @(6, 7)
    7  = @(6, 6)
    8  = @(1, 3)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(1, 6)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 5)
    15 = @(1, 5)
    16 = @(1, 4)
    17 = @(1, 4)
    18 = @(1, 5)
    19 = @(1, 4)
    20 = @(1, 3)
    21 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
